$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.272677421569824
$ws.Range("B1").Value = 3.229651927947998
$ws.Range("C1").Value = 4.005142211914062
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 3.198182821273804
